# "end slide + export to ppt, pdf"
#
# Adds a new final "The end" slide (slide 17) to the deck. The new slide
# reuses the same placeholder layout ("Titelfolie" / ctrTitle + footer +
# date + slide-number) as slide 10 ("Demo"), so the cleanest way to get an
# identical shape/placeholder structure is to duplicate slide 10 and move
# the duplicate to the end of the deck, then edit its text.

$p = $ppt.ActivePresentation

# --- duplicate the slide that already has the matching placeholder set ---
$srcSlide = $p.Slides.Item(10)
$dupRange = $srcSlide.Duplicate()
$newSlide = $dupRange.Item(1)

# move the duplicate to the very end of the deck (position 17)
$newSlide.MoveTo($p.Slides.Count)

# --- title placeholder: "The end" ---
$title = $newSlide.Shapes.Item(1)
$tr = $title.TextFrame.TextRange
$tr.Text = "The"
$tr = $tr.InsertAfter(" ")
$tr = $tr.InsertAfter("end")

# --- footer placeholder: "Master Lab Course Web Applications" ---
$footer = $newSlide.Shapes.Item(2)
$ftr = $footer.TextFrame.TextRange
$ftr.Text = "Master Lab "
$ftr = $ftr.InsertAfter("Course")
$ftr = $ftr.InsertAfter(" Web ")
$ftr = $ftr.InsertAfter("Applications")

# footer placeholder also got an explicit (resized/centered) position
$footer.Left = 191.2488188976378
$footer.Top = 504
$footer.Width = 490.7511811023622
$footer.Height = 24

# date + slide-number placeholders (Shapes.Item(3)/(4)) keep their
# inherited "25.02.2013" / auto slide-number field content unchanged.
